$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old numeric cells A1, B1, C1
$ws.Range("A1:C1").ClearContents()

# Put the new text value into B2
$ws.Range("B2").Value = "5+3+4=12"

# Update selection to match the new active cell
$ws.Range("B2").Select()
